$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 2 (shifts existing data rows 2-21 down to 4-23)
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()
$ws.Range("A2:C3").ClearFormats()

# New data for the two inserted rows
$newTopRows = @(
    @(0.5199990272521973, -4.081954479217529, 1.088714599609375),
    @(0.446847915649414, -3.855323314666748, 1.094059705734253)
)

for ($i = 0; $i -lt $newTopRows.Length; $i++) {
    $r = 2 + $i
    $row = $newTopRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# New rows to append at the bottom (rows 24-31)
$newBottomRows = @(
    @(-0.7247915863990784, -2.964529037475586, -2.036930084228516),
    @(-0.9755517840385436, -3.013092756271362, -0.1954768747091293),
    @(-0.0148134818300604, -4.330729007720947, 0.6641632318496704),
    @(-0.2535090744495392, -4.50192403793335, 0.8677340745925903),
    @(-0.0308486949652433, -3.680310487747192, 1.009607553482056),
    @(-0.2924517393112182, 0.6568328738212585, 0.4216497242450714),
    @(0.683863639831543, 4.383111000061035, -1.505782842636108),
    @(0.6068946123123169, 4.862334728240967, -0.4990769028663635)
)

for ($i = 0; $i -lt $newBottomRows.Length; $i++) {
    $r = 24 + $i
    $row = $newBottomRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
